# Update betting odds values in Sheet1 as per the 2024-11-28 FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.25
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 2.1
$ws.Range("L2").Value = 3.6
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 11
$ws.Range("O2").Value = 1.3
$ws.Range("P2").Value = 3.4
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 1.8
$ws.Range("S2").Value = 1.4
$ws.Range("T2").Value = 2.75
$ws.Range("U2").Value = 1.73
$ws.Range("V2").Value = 2
$ws.Range("W2").Value = 8
$ws.Range("Y2").Value = 9.5
$ws.Range("Z2").Value = 21
$ws.Range("AA2").Value = 19
$ws.Range("AB2").Value = 29
$ws.Range("AC2").Value = 10
$ws.Range("AD2").Value = 6.5
$ws.Range("AE2").Value = 15
$ws.Range("AF2").Value = 51
$ws.Range("AG2").Value = 201
$ws.Range("AH2").Value = 9.5
$ws.Range("AI2").Value = 15
$ws.Range("AJ2").Value = 11
$ws.Range("AM2").Value = 34
$ws.Range("AO2").Value = 13
$ws.Range("AP2").Value = 23
$ws.Range("AR2").Value = 67
$ws.Range("AS2").Value = 151
$ws.Range("AT2").Value = 2.75
$ws.Range("AU2").Value = 8
$ws.Range("AW2").Value = 5
$ws.Range("AY2").Value = 26
$ws.Range("BA2").Value = 81

# Row 3
$ws.Range("Q3").Value = 1.7
$ws.Range("R3").Value = 2.1

# Row 6
$ws.Range("G6").Value = 2.2
$ws.Range("J6").Value = 3
$ws.Range("L6").Value = 3.6
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 3
$ws.Range("X6").Value = 11
$ws.Range("Y6").Value = 10
$ws.Range("AA6").Value = 21
$ws.Range("AB6").Value = 34
$ws.Range("AC6").Value = 8.5
$ws.Range("AX6").Value = 17

# Row 7
$ws.Range("G7").Value = 1.85
$ws.Range("I7").Value = 3.75
$ws.Range("L7").Value = 4.33
$ws.Range("O7").Value = 1.3
$ws.Range("P7").Value = 3.4
$ws.Range("Q7").Value = 1.98
$ws.Range("R7").Value = 1.83
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.83
$ws.Range("X7").Value = 9
$ws.Range("Z7").Value = 15
$ws.Range("AG7").Value = 700
$ws.Range("AH7").Value = 11
$ws.Range("AI7").Value = 21
$ws.Range("AL7").Value = 34
$ws.Range("AM7").Value = 41
$ws.Range("AW7").Value = 6
$ws.Range("BA7").Value = 101

# Row 9
$ws.Range("G9").Value = 2.7
$ws.Range("I9").Value = 2.4
$ws.Range("L9").Value = 3.25
$ws.Range("Q9").Value = 2.35
$ws.Range("R9").Value = 1.57
$ws.Range("X9").Value = 13
$ws.Range("AI9").Value = 11
$ws.Range("AK9").Value = 23
